# Applies the "Bao cao, de cuong 31/8" edit to the presentation.
# Units: Shape Left/Top/Width/Height are in points (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape "Flowchart: Magnetic Disk 26" (id=27): move up & shrink height
# off 0,4038600 -> 0,3962400 ; ext 1524000x1828800 -> 1524000x1524000
$diskShape = $s1.Shapes.Item("Flowchart: Magnetic Disk 26")
$diskShape.Left = 0
$diskShape.Top = 3962400 / 12700
$diskShape.Width = 1524000 / 12700
$diskShape.Height = 1524000 / 12700

# Shape "TextBox 32" ("Title"): shift right
# off 1600200,3962400 -> 1752600,3962400
$titleBox = $s1.Shapes.Item("TextBox 32")
$titleBox.Left = 1752600 / 12700
$titleBox.Top = 3962400 / 12700

# Shape "TextBox 33" ("Title, Abtrach" -> "Title, abstract")
$titleAbstractBox = $s1.Shapes.Item("TextBox 33")
$tr = $titleAbstractBox.TextFrame.TextRange
$run2 = $tr.Characters(8, 7)
$run2.Text = "abstract"

# Delete shapes "Right Arrow 34" (id=35) and "TextBox 35" (id=36)
$s1.Shapes.Item("Right Arrow 34").Delete()
$s1.Shapes.Item("TextBox 35").Delete()

# ---------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Shape "TextBox 52" ("Rut" -> "Lay thong Tin"), reposition/resize and reword
$rutBox = $s4.Shapes.Item("TextBox 52")
$trR = $rutBox.TextFrame.TextRange
$trR.Text = "Lấy thông Tin"
$trR.Font.Size = 14
$rR1 = $trR.Characters(1, 3)
$rR1.Text = "Lấy"
$rR2 = $trR.Characters(4, 1)
$rR2.Text = " "
$rR3 = $trR.Characters(5, 5)
$rR3.Text = "thông"
$rR4 = $trR.Characters(10, 4)
$rR4.Text = " Tin"

$rutBox.Left = 1600200 / 12700
$rutBox.Top = 2590800 / 12700
$rutBox.Width = 1752600 / 12700
$rutBox.Height = 24.23445

# Shape "TextBox 53" ("(3) " + "Module " + "Import " -> merged "(3) Module Import ")
$moduleBox = $s4.Shapes.Item("TextBox 53")
$trM = $moduleBox.TextFrame.TextRange
$mergedRun = $trM.Characters(1, 18)
$mergedRun.Text = "(3) Module Import "
